function Set-EventRow {
    param($ws, $row, $b, $c, $d, $e, $f, $g, $h, $i)
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).NumberFormat = "General"
    $ws.Cells.Item($row, 6).Value = $f
    if ($g -is [string]) {
        $ws.Cells.Item($row, 7).NumberFormat = "@"
    } else {
        $ws.Cells.Item($row, 7).NumberFormat = "General"
    }
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
    $ws.Cells.Item($row, 9).Value = $i
}

function Update-EventSheet {
    param($ws)

    # Force column B (start-date strings like "2024-03-23") to be stored as
    # plain text, not auto-converted to Excel date serials.
    $ws.Range("B2:B23").NumberFormat = "@"

    Set-EventRow $ws 2 "2024-03-23" "上饶·原×铁×崩only（取消）" "五三东大道42号 回禾酒店" "2024.03.23 10:00-03.23 17:00" 44 "不可售" "https://show.bilibili.com/platform/detail.html?id=81103" "//i2.hdslb.com/bfs/openplatform/202401/pp6c5TsC1705647180602.jpeg"
    Set-EventRow $ws 3 "2024-03-23" "南昌·AP动漫游戏嘉年华" "八一桥街道青山南路118号 蓝海会展中心" "2024.03.23 09:00-03.24 17:00" 2891 58.5 "https://show.bilibili.com/platform/detail.html?id=81232" "//i2.hdslb.com/bfs/openplatform/202401/NZv97SmS1705912230957.jpeg"
    Set-EventRow $ws 4 "2024-03-23" "南昌·运动番only春季集训（取消）" "创新三路777号 南昌小飞侠章鱼文化体育公园" "2024.03.23 10:00-03.24 17:00" 189 "不可售" "https://show.bilibili.com/platform/detail.html?id=81950" "//i1.hdslb.com/bfs/openplatform/202402/bm4uH4qB1708425538357.jpeg"
    Set-EventRow $ws 5 "2024-03-24" "南昌·AP动漫游戏  嘉年华内场票-小N&子音" "八一桥街道青山南路118号 蓝海会展中心" "2024.03.24 09:00-03.24 17:00" 148 218 "https://show.bilibili.com/platform/detail.html?id=81973" "//i0.hdslb.com/bfs/openplatform/202402/zbG5HICL1708504962467.jpeg"
    Set-EventRow $ws 6 "2024-03-24" "鹰潭·宅舞联萌·随舞动漫派对（免费活动)" "玉清路与象山路交叉口东南角 鹰潭天虹购物中心" "2024.03.24 14:00-03.24 18:00" 13 22.33 "https://show.bilibili.com/platform/detail.html?id=82434" "//i0.hdslb.com/bfs/openplatform/202403/oj5AEi4W1709713367717.jpeg"
    Set-EventRow $ws 7 "2024-03-30" "南昌·CM01动漫游戏博览会" "怀玉山大道1315号 南昌绿地国际博览中心" "2024.03.30 10:00-03.31 17:00" 1551 55 "https://show.bilibili.com/platform/detail.html?id=81691" "//i0.hdslb.com/bfs/openplatform/202403/mKGcIL5D1710500376094.png"
    Set-EventRow $ws 8 "2024-03-30" "鹰潭·原×铁×崩only（取消）" "南站路24号 回禾酒店(鹰潭火车站店)" "2024.03.30 10:00-03.30 17:00" 36 "不可售" "https://show.bilibili.com/platform/detail.html?id=81097" "//i2.hdslb.com/bfs/openplatform/202401/q0AZaXAk1705646244207.jpeg"
    Set-EventRow $ws 9 "2024-03-31" "新余·文旅国漫嘉年华暨BM次元盛典" "五一南路与仙女湖大道交叉口西北 老上海风情街白金汉宫" "2024.03.31 10:00-03.31 17:00" 76 60 "https://show.bilibili.com/platform/detail.html?id=82208" "//i2.hdslb.com/bfs/openplatform/202403/aXc6vPDP1709547191851.jpeg"
    Set-EventRow $ws 10 "2024-03-31" "景德镇·宅舞联萌·随舞动漫派对（免费活动)" "经二路与纬二路交叉路口 景德镇市宝龙广场" "2024.03.31 14:00-03.31 18:00" 23 22.33 "https://show.bilibili.com/platform/detail.html?id=82437" "//i0.hdslb.com/bfs/openplatform/202403/zcNNsicc1709714035066.jpeg"
    Set-EventRow $ws 11 "2024-04-04" "南昌·创造力动漫游戏嘉年华1.0" "八一桥街道青山南路118号 蓝海会展中心" "2024.04.04 10:00-04.05 17:00" 1282 55 "https://show.bilibili.com/platform/detail.html?id=82419" "//i2.hdslb.com/bfs/openplatform/202403/vNXlqMyb1710471001741.jpeg"
    Set-EventRow $ws 12 "2024-04-04" "南昌·梦想家动漫展" "龙蟠街666号融创茂1层 融创茂" "2024.04.04 10:00-04.04 17:00" 9 45 "https://show.bilibili.com/platform/detail.html?id=82756" "//i0.hdslb.com/bfs/openplatform/202403/zPhYpquw1710235466897.png"
    Set-EventRow $ws 13 "2024-04-04" "赣州·第三届半夏动漫展" "105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心" "2024.04.04 10:00-04.06 17:00" 411 50 "https://show.bilibili.com/platform/detail.html?id=82235" "//i0.hdslb.com/bfs/openplatform/202403/4DWZWYGm1709278879159.jpeg"
    Set-EventRow $ws 14 "2024-04-04" "赣州·赣次元·归来国风动漫节" "客家大道568号文清外国语学校旁 赣州市文清外国语学校国际交流中心" "2024.04.04 10:00-04.04 17:00" 334 40 "https://show.bilibili.com/platform/detail.html?id=82125" "//i1.hdslb.com/bfs/openplatform/202403/fIehikk51709705287036.jpeg"
    Set-EventRow $ws 15 "2024-04-05" "抚州·第七届FZ动漫文化节" "迎宾大道288号 凤凰世纪名都大酒店" "2024.04.05 09:30-04.05 17:00" 62 50 "https://show.bilibili.com/platform/detail.html?id=82381" "//i2.hdslb.com/bfs/openplatform/202403/Y725SN0G1709694367526.jpeg"
    Set-EventRow $ws 16 "2024-04-06" "萍乡·2024DDS国漫盛典" "凤凰街迎宾路18号 鸿凯大酒店" "2024.04.06 10:00-04.06 17:00" 51 50 "https://show.bilibili.com/platform/detail.html?id=82413" "//i1.hdslb.com/bfs/openplatform/202403/Rbu1xyFr1709707847098.jpeg"
    Set-EventRow $ws 17 "2024-04-13" "南昌·原X穹X崩only" "丰和北大道299号 新吉花园酒店" "2024.04.13 10:00-04.13 17:00" 114 65 "https://show.bilibili.com/platform/detail.html?id=80807" "//i0.hdslb.com/bfs/openplatform/202402/kfK13XvH1709202705153.jpeg"
    Set-EventRow $ws 18 "2024-04-13" "南昌·第二届漫拥动漫嘉年华mini" "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆" "2024.04.13 10:00-04.14 18:00" 83 39.9 "https://show.bilibili.com/platform/detail.html?id=82210" "//i0.hdslb.com/bfs/openplatform/202402/KYd0bfk11709203777701.png"
    Set-EventRow $ws 19 "2024-04-20" "南昌·DSL国风动漫游戏嘉年华" "九龙大道1388号（上饶街与九龙大道交叉口西北100米） 中国南昌虚拟现实VR产业基地" "2024.04.20 09:00-04.21 17:00" 95 35 "https://show.bilibili.com/platform/detail.html?id=82107" "//i0.hdslb.com/bfs/openplatform/202402/QDlumVb41708943318282.jpeg"
    Set-EventRow $ws 20 "2024-04-20" "南昌·New World国潮动漫博览会" "怀玉山大道1315号 南昌绿地国际博览中心" "2024.04.20 09:30-04.21 17:00" 2903 60 "https://show.bilibili.com/platform/detail.html?id=82411" "//i1.hdslb.com/bfs/openplatform/202403/xbYbLXc81709707724935.jpeg"
    Set-EventRow $ws 21 "2024-05-01" "九江·第三届ACD动漫游戏嘉年华" "九瑞大道与重庆路交汇处西南角 九江国际会展中心" "2024.05.01 09:00-05.02 17:00" 350 55 "https://show.bilibili.com/platform/detail.html?id=82464" "//i0.hdslb.com/bfs/openplatform/202403/HjMMyP3a1709780146797.jpeg"
    Set-EventRow $ws 22 "2024-05-02" "抚州·临次元07国漫&运动番嘉年华" "伍塘路1098号 乐课篮球公园" "2024.05.02 10:00-05.02 16:00" 13 40 "https://show.bilibili.com/platform/detail.html?id=82911" "//i0.hdslb.com/bfs/openplatform/202403/Pc5TMZ001710429899787.jpeg"
    Set-EventRow $ws 23 "2024-05-26" "南昌·代号鸢盛花行only" "民德路411号 东方豪景花园酒店(民德路店)" "2024.05.26 09:30-05.26 17:30" 68 78 "https://show.bilibili.com/platform/detail.html?id=82529" "//i1.hdslb.com/bfs/openplatform/202403/TJ8iC73c1709804909450.png"

    # The last existing row (24) has no successor event anymore once the
    # top entry is dropped and everything above moves up one slot, so the
    # sheet now ends at row 23 (A1:I23) instead of row 24 (A1:I24).
    $ws.Rows.Item(24).Delete()
}

$wb = $excel.ActiveWorkbook

# Both the "展览" (sheet 1) and "全部类型" (sheet 4) tabs carry the same
# event table and both need the refreshed data.
Update-EventSheet $wb.Worksheets.Item("展览")
Update-EventSheet $wb.Worksheets.Item("全部类型")

